$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column G should become "Tumakuru (Tumkur)"
$tumakuruRows = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,41,42,43)

foreach ($r in $tumakuruRows) {
    $ws.Cells.Item($r, 7).Value = "Tumakuru (Tumkur)"
}

# Row 20: Yadagiri -> Yadgir
$ws.Cells.Item(20, 7).Value = "Yadgir"

# Row 23: remove the empty F23 cell entirely (clear its contents)
$ws.Cells.Item(23, 6).ClearContents()
